$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Zelem"
$ws.Range("B2").Value = 82
$ws.Range("E2").Value = "CDM"
$ws.Range("H2").Value = "Barclays WSL"
$ws.Range("N2").Value = 750
$ws.Range("P2").Value = 147

# Row 3
$ws.Range("A3").Value = "van de Donk"
$ws.Range("B3").Value = 82
$ws.Range("E3").Value = "CAM"
$ws.Range("H3").Value = "D1 Arkema"
$ws.Range("I3").Value = "OL"
$ws.Range("N3").Value = 750
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 169

# Row 4
$ws.Range("A4").Value = "Lopes"
$ws.Range("D4").Value = "Rare"
$ws.Range("E4").Value = "GK"
$ws.Range("G4").Value = "Portugal"
$ws.Range("H4").Value = "Ligue 1 Uber Eats"
$ws.Range("I4").Value = "OL"
$ws.Range("N4").Value = 750
$ws.Range("P4").Value = 171

# Row 5
$ws.Range("A5").Value = "Pacheco"
$ws.Range("B5").Value = 81
$ws.Range("D5").Value = "Rare"
$ws.Range("E5").Value = "LB"
$ws.Range("H5").Value = "Barclays WSL"
$ws.Range("I5").Value = "Aston Villa"
$ws.Range("N5").Value = 700
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 205

# Row 6
$ws.Range("A6").Value = "Le Tissier"
$ws.Range("B6").Value = 81
$ws.Range("E6").Value = "CB"
$ws.Range("G6").Value = "England"
$ws.Range("H6").Value = "Barclays WSL"
$ws.Range("I6").Value = "Manchester Utd"
$ws.Range("P6").Value = 226

# Row 7
$ws.Range("A7").Value = "Ramos"
$ws.Range("B7").Value = 80
$ws.Range("D7").Value = "Rare"
$ws.Range("E7").Value = "ST"
$ws.Range("G7").Value = "Portugal"
$ws.Range("I7").Value = "Paris SG"
$ws.Range("N7").Value = 700
$ws.Range("O7").Value = 3
$ws.Range("P7").Value = 296

# Row 8
$ws.Range("A8").Value = "Groenen"
$ws.Range("B8").Value = 79
$ws.Range("E8").Value = "CDM"
$ws.Range("H8").Value = "D1 Arkema"
$ws.Range("I8").Value = "Paris SG"
$ws.Range("P8").Value = 347

# Row 9
$ws.Range("A9").Value = "Cissoko"
$ws.Range("B9").Value = 79
$ws.Range("E9").Value = "CB"
$ws.Range("G9").Value = "France"
$ws.Range("H9").Value = "Barclays WSL"
$ws.Range("I9").Value = "West Ham"
$ws.Range("N9").Value = 500
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 372

# Row 10
$ws.Range("A10").Value = "Cambot"
$ws.Range("B10").Value = 77
$ws.Range("E10").Value = "ST"
$ws.Range("G10").Value = "France"
$ws.Range("H10").Value = "D1 Arkema"
$ws.Range("I10").Value = "En Avant Guingamp"
$ws.Range("N10").Value = 500
$ws.Range("P10").Value = 497

# Row 11
$ws.Range("A11").Value = "Guilbert"
$ws.Range("B11").Value = 76
$ws.Range("D11").Value = "Common"
$ws.Range("G11").Value = "France"
$ws.Range("I11").Value = "Strasbourg"
$ws.Range("N11").Value = 500
$ws.Range("P11").Value = 635

# Row 12
$ws.Range("A12").Value = "Jurić"
$ws.Range("B12").Value = 59
$ws.Range("C12").Value = "Bronze"
$ws.Range("E12").Value = "ST"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "Bosnia Herzegovina"
$ws.Range("H12").Value = "PKO BP Ekstraklasa"
$ws.Range("I12").Value = "ŁKS Łódź"
$ws.Range("N12").Value = 200
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 833
